$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D; this shifts the existing D:K quarterly
# data (and its formatting) two columns to the right, into F:M.
$ws.Columns("D:E").Insert()

# Seed the two new columns (D:E) with formatting copied from the (now shifted)
# old D:E data, which now lives in F:G, for every row that has data.
$ws.Range("F7:G35").Copy($ws.Range("D7:E35"))
$ws.Range("F38:G77").Copy($ws.Range("D38:E77"))
$ws.Range("F80:G102").Copy($ws.Range("D80:E102"))

# Write the two new quarters of data into D:E for each line item.
$ws.Range("D7").Value = 43434
$ws.Range("E7").Value = 43343
$ws.Range("D8").Value = 112000
$ws.Range("E8").Value = 123300
$ws.Range("D9").Value = 83300
$ws.Range("E9").Value = 91000
$ws.Range("D10").Value = 28700
$ws.Range("E10").Value = 32300
$ws.Range("D12").Value = 3600
$ws.Range("E12").Value = 4100
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 109900
$ws.Range("E17").Value = 116700
$ws.Range("D18").Value = 2100
$ws.Range("E18").Value = 6600
$ws.Range("D20").Value = 800
$ws.Range("E20").Value = 500
$ws.Range("D21").Value = 6300
$ws.Range("E21").Value = 10800
$ws.Range("D22").Value = 1200
$ws.Range("E22").Value = 1200
$ws.Range("D23").Value = 1700
$ws.Range("E23").Value = 5900
$ws.Range("D24").Value = 500
$ws.Range("E24").Value = 1500
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 1200
$ws.Range("E26").Value = 4500
$ws.Range("D27").Value = 1200
$ws.Range("E27").Value = 4500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = 500
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -800
$ws.Range("E32").Value = -500
$ws.Range("D33").Value = 1200
$ws.Range("E33").Value = 5000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 1200
$ws.Range("E35").Value = 5000
$ws.Range("D38").Value = 43434
$ws.Range("E38").Value = 43343
$ws.Range("D41").Value = 137200
$ws.Range("E41").Value = 160800
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 84900
$ws.Range("E43").Value = 69100
$ws.Range("D44").Value = 88900
$ws.Range("E44").Value = 79200
$ws.Range("D45").Value = 14300
$ws.Range("E45").Value = 21900
$ws.Range("D46").Value = 325300
$ws.Range("E46").Value = 331100
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 60500
$ws.Range("E48").Value = 57200
$ws.Range("D49").Value = 91100
$ws.Range("E49").Value = 92000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 25600
$ws.Range("E52").Value = 19900
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 502500
$ws.Range("E54").Value = 500300
$ws.Range("D57").Value = 41300
$ws.Range("E57").Value = 30500
$ws.Range("D58").Value = 200
$ws.Range("E58").Value = 200
$ws.Range("D59").Value = 41500
$ws.Range("E59").Value = 49400
$ws.Range("D60").Value = 83000
$ws.Range("E60").Value = 80100
$ws.Range("D61").Value = 116500
$ws.Range("E61").Value = 116600
$ws.Range("D62").Value = 26300
$ws.Range("E62").Value = 26700
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 225800
$ws.Range("E66").Value = 223400
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 483800
$ws.Range("E72").Value = 484900
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 276700
$ws.Range("E76").Value = 276900
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43434
$ws.Range("E80").Value = 43343
$ws.Range("D81").Value = 1200
$ws.Range("E81").Value = 5000
$ws.Range("D83").Value = 3400
$ws.Range("E83").Value = 3700
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -14300
$ws.Range("E89").Value = 25800
$ws.Range("D91").Value = -5700
$ws.Range("E91").Value = -4100
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -4700
$ws.Range("E94").Value = 27800
$ws.Range("D96").Value = -3300
$ws.Range("E96").Value = -3300
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -4300
$ws.Range("E100").Value = -3400
$ws.Range("D101").Value = -200
$ws.Range("E101").Value = -1100
$ws.Range("D102").Value = -23600
$ws.Range("E102").Value = 49000
